# memperbaiki logika untuk menyimpan multiple data dan memperbaiki tabel angsuran
#
# The "pinjaman mandiri / kredit / pinjaman" sub-category row (row 12) was a
# duplicate/incorrect entry for the installment ("angsuran") table, so remove
# it entirely. Excel shifts the remaining rows (pembayaran angsuran, simpanan
# pokok) up by one and prunes the now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing "pinjaman mandiri" / "kredit" / "pinjaman"
$ws.Rows.Item(12).Delete()

# Restore the cursor/selection position left by the author's last save
$ws.Range("J3").Select()
